$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.916.60"
$ws.Range("E2").Value = "  -0.81%  "
$ws.Range("D3").Value = "2.432.25"
$ws.Range("E3").Value = "  -0.49%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "553.77"
$ws.Range("E5").Value = "  -0.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.72"
$ws.Range("E6").Value = "  -0.42%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.510"
$ws.Range("E8").Value = "  +1.96%  "
$ws.Range("E9").Value = "  +8.16%  "
$ws.Range("E10").Value = "  -0.59%  "
$ws.Range("E11").Value = "  -0.95%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.78"
$ws.Range("E12").Value = "  +0.31%  "
$ws.Range("D13").Value = "67.833.97"
$ws.Range("E13").Value = "  -0.55%  "
$ws.Range("E14").Value = "  +1.67%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "22.98"
$ws.Range("E15").Value = "  -0.92%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "10.35"
$ws.Range("E16").Value = "  -2.97%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "333.83"
$ws.Range("E17").Value = "  -1.34%  "
$ws.Range("E18").Value = "  -1.83%  "
$ws.Range("E19").Value = "  +0.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.85"
$ws.Range("E21").Value = "  +0.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "66.12"
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.63"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.09"
$ws.Range("E24").Value = "  +0.70%  "
$ws.Range("D25").Value = "0.0₃0808"
$ws.Range("E25").Value = "  +0.22%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.09"
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "417.50"
$ws.Range("E28").Value = "  -3.21%  "
$ws.Range("E29").Value = "  +1.12%  "
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "160.79"
$ws.Range("E31").Value = "  +2.99%  "
$ws.Range("E32").Value = "  -0.42%  "
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.84"
$ws.Range("E34").Value = "  +0.59%  "
$ws.Range("E35").Value = "  -3.49%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.294"
$ws.Range("E36").Value = "  -2.08%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.25"
$ws.Range("E37").Value = "  -2.71%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.46"
$ws.Range("E38").Value = "  +1.25%  "
$ws.Range("E39").Value = "  -0.67%  "
$ws.Range("E40").Value = "  -1.03%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.32"
$ws.Range("E41").Value = "  +0.33%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "129.54"
$ws.Range("E42").Value = "  -1.29%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0707"
$ws.Range("E43").Value = "  -0.61%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.478"
$ws.Range("E44").Value = "  -0.16%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.555"
$ws.Range("E45").Value = "  -0.25%  "
$ws.Range("E46").Value = "  +1.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.33"
$ws.Range("E48").Value = "  -6.38%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "16.53"
$ws.Range("E49").Value = "  -0.86%  "
$ws.Range("E50").Value = "  +4.01%  "
$ws.Range("E51").Value = "  +0.79%  "
